$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pending-approval rows (20-25) mirror the layout/style of row 19
# (same party/company/bank details), so copy row 19 as a template for
# formatting + the shared columns, then overwrite the per-row fields.
for ($row = 20; $row -le 25; $row++) {
    $ws.Range("A19:AO19").Copy($ws.Range("A" + $row + ":AO" + $row))
}

$rows = @(
    @{ Row = 20; UniqueId = "88538b2f-9196-4e71-8def-4bbb1d589392"; Amount = 7865;   Narration = "Train ticket to Goa RPA_UNIQUE_ID : 7b63b7fd-265c-4124-a7a5-5d8044a3ddf1"; Project = "IPSHEM"; },
    @{ Row = 21; UniqueId = "b6315e8f-4f01-486f-a315-f31da08cbb3b"; Amount = 239350; Narration = "194C RPA_UNIQUE_ID : 7acec7df-0570-4edb-a065-bb1cdcb9c474"; Project = "IPSHEM"; },
    @{ Row = 22; UniqueId = "df7c0414-0a8f-4271-8fc1-8fb47bcbb273"; Amount = 20000;  Narration = "194C RPA_UNIQUE_ID : e1a23f47-3131-451b-bc65-27b4876a3d21"; Project = "Cochin"; },
    @{ Row = 23; UniqueId = "0f2b8be4-bfbb-46e3-af80-74d6df7d2161"; Amount = 163992; Narration = "194A RPA_UNIQUE_ID : d3c85736-1f30-421a-9331-c7b682a12595"; Project = "IPSHEM"; },
    @{ Row = 24; UniqueId = "7de1fff2-50f3-4fd3-96eb-2945e899e145"; Amount = 9500;   Narration = "194I RPA_UNIQUE_ID : 5e874df5-7e3c-4aff-81b4-c5e8395eeda0"; Project = "IPSHEM"; },
    @{ Row = 25; UniqueId = "6d899d84-a48e-426f-ad38-3f931a1504d6"; Amount = 16540;  Narration = "194J RPA_UNIQUE_ID : 48b65542-4c14-407f-85e4-297c58de7d36"; Project = "Cochin"; }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 12).Value = $r.UniqueId                 # L - UNIQUE_ID
    $ws.Cells.Item($row, 21).Value = "pending"                   # U - STATUS
    $ws.Cells.Item($row, 22).Value = $r.Amount                   # V - BASIC_AMOUNT
    $ws.Cells.Item($row, 24).Value = $r.Narration                # X - NARRATION
    $ws.Cells.Item($row, 25).Value = $r.Project                  # Y - PROJECT_NAME
    $ws.Cells.Item($row, 26).Value = "PAYMENT"                   # Z - CATEGORY
    $ws.Cells.Item($row, 27).Value = "Payments@westernidc.com"   # AA - FROM_MAIL
    $ws.Cells.Item($row, 28).Value = "ESTIMATION NOT MATCHED"    # AB - STATUS_MATCHED_ESTIMATION
    $ws.Cells.Item($row, 29).Value = 0                           # AC - FIXED_AMOUNT
    $ws.Cells.Item($row, 30).Value = 0                           # AD - BALANCE_AMOUNT
    $ws.Cells.Item($row, 31).Value = 0                           # AE - ADJUSTMENT_AMOUNT
}
